$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A21").Value = "Chart with  realtime data"
$ws.Range("B21").Value = "Charts"
$ws.Range("C21").Value = "Displays a chart that continuously refreshes with updated realtime date.  This example uses random data generated in the Edit Chart > Script >Interactivity > Chart Area > Load."
$ws.Range("D21").Value = "Clement Wong"
$ws.Range("E21").Value = "#experts-BIRT email sent 27/23/2014 2:46 PM Eastern"
$ws.Range("F21").Value = "Report Designs/Charts/HTML5 Charts - Dynamic Updating Highchart__cwong.rptdesign"

$ws.Range("F22").Select()
